# Updates from Issues 42/43/18
# Mark Genesis family (rows 15-18) and PSX (row 25) Functional Test column as "Complete"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D15").Value = "Complete"
$ws.Range("D16").Value = "Complete"
$ws.Range("D17").Value = "Complete"
$ws.Range("D18").Value = "Complete"
$ws.Range("D25").Value = "Complete"

# Reflect the final selection state (whole used range, active cell D33)
$ws.Range("A1:D33").Select() | Out-Null
